$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 198
$ws.Range("I6").Value = 45
$ws.Range("K6").Value = 135
$ws.Range("M6").Value = -23
$ws.Range("H17").Value = 9952.35
$ws.Range("J17").Value = 10376.211
$ws.Range("L17").Value = 31128.633
$ws.Range("N17").Value = -31464.633
$ws.Range("H64").Value = 2923.1904
$ws.Range("I64").Value = 2737.6
$ws.Range("J64").Value = 2981.1875
$ws.Range("K64").Value = 2737.6
$ws.Range("L64").Value = 2981.1875
$ws.Range("M64").Value = -2489.6
$ws.Range("N64").Value = -3477.1875
$ws.Range("H67").Value = 2923.1904
$ws.Range("I67").Value = 2737.6
$ws.Range("J67").Value = 2981.1875
$ws.Range("K67").Value = 2737.6
$ws.Range("L67").Value = 2981.1875
$ws.Range("M67").Value = -1879.6
$ws.Range("N67").Value = -4697.1875
$ws.Range("H92").Value = 38461924
$ws.Range("I92").Value = 55555790
$ws.Range("J92").Value = 723.125
$ws.Range("K92").Value = 55555790
$ws.Range("L92").Value = 723.125
$ws.Range("M92").Value = -55554542
$ws.Range("N92").Value = -3219.125
$ws.Range("H96").Value = 796.2857
$ws.Range("I96").Value = 832
$ws.Range("J96").Value = 769.5
$ws.Range("K96").Value = 2496
$ws.Range("L96").Value = 2308.5
$ws.Range("M96").Value = -1123
$ws.Range("N96").Value = -5054.5
$ws.Range("H112").Value = 2925126.2
$ws.Range("J112").Value = 2925126.2
$ws.Range("L112").Value = 8775378.600000001
$ws.Range("N112").Value = -8777594.600000001
$ws.Range("H127").Value = 1247.7059
$ws.Range("I127").Value = 928.1429
$ws.Range("J127").Value = 1471.4
$ws.Range("K127").Value = 2784.4287
$ws.Range("L127").Value = 4414.200000000001
$ws.Range("M127").Value = 2175.5713
$ws.Range("N127").Value = -14334.2
$ws.Range("H129").Value = 1107.0454
$ws.Range("J129").Value = 1224.4736
$ws.Range("L129").Value = 3673.4208
$ws.Range("N129").Value = -13673.4208
$ws.Range("H138").Value = 4310.4165
$ws.Range("I138").Value = 3686.375
$ws.Range("J138").Value = 4622.4375
$ws.Range("K138").Value = 11059.125
$ws.Range("L138").Value = 13867.3125
$ws.Range("M138").Value = -5919.125
$ws.Range("N138").Value = -24147.3125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5082.9214
$ws.Range("I32").Value = 4004.6956
$ws.Range("K32").Value = 4004.6956
$ws.Range("M32").Value = -3717.6956
$ws.Range("H61").Value = 3050.4062
$ws.Range("I61").Value = 2820.4333
$ws.Range("J61").Value = 6500
$ws.Range("K61").Value = 2820.4333
$ws.Range("L61").Value = 6500
$ws.Range("M61").Value = -2608.4333
$ws.Range("N61").Value = -6924
$ws.Range("H74").Value = 100004110
$ws.Range("I74").Value = 166669680
$ws.Range("J74").Value = 5749.75
$ws.Range("K74").Value = 166669680
$ws.Range("L74").Value = 5749.75
$ws.Range("M74").Value = -166668806
$ws.Range("N74").Value = -7497.75
$ws.Range("H77").Value = 100004110
$ws.Range("I77").Value = 166669680
$ws.Range("J77").Value = 5749.75
$ws.Range("K77").Value = 833348400
$ws.Range("L77").Value = 28748.75
$ws.Range("M77").Value = -833344032
$ws.Range("N77").Value = -37484.75
$ws.Range("H102").Value = 856.6111
$ws.Range("I102").Value = 856.6111
$ws.Range("K102").Value = 856.6111
$ws.Range("M102").Value = 765.3889
$ws.Range("H122").Value = 2552.75
$ws.Range("I122").Value = 2181.4443
$ws.Range("J122").Value = 3221.1
$ws.Range("K122").Value = 6544.3329
$ws.Range("L122").Value = 9663.3
$ws.Range("M122").Value = -4094.3329
$ws.Range("N122").Value = -14563.3
$ws.Range("H136").Value = 3050.4062
$ws.Range("I136").Value = 2820.4333
$ws.Range("J136").Value = 6500
$ws.Range("K136").Value = 8461.2999
$ws.Range("L136").Value = 19500
$ws.Range("M136").Value = -5911.2999
$ws.Range("N136").Value = -24600

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 896.6667
$ws.Range("I107").Value = 821.25
$ws.Range("K107").Value = 821.25
$ws.Range("M107").Value = 1098.75
$ws.Range("H134").Value = 3686.1875
$ws.Range("I134").Value = 4146.1113
$ws.Range("J134").Value = 1202.6
$ws.Range("K134").Value = 12438.3339
$ws.Range("L134").Value = 3607.8
$ws.Range("M134").Value = -9903.333899999998
$ws.Range("N134").Value = -8677.8
$ws.Range("H137").Value = 50640
$ws.Range("J137").Value = 50640
$ws.Range("L137").Value = 50640
$ws.Range("N137").Value = -60840

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 554.3333
$ws.Range("J22").Value = 751
$ws.Range("L22").Value = 751
$ws.Range("N22").Value = -1451
$ws.Range("H31").Value = 3002.5334
$ws.Range("I31").Value = 2401.8125
$ws.Range("J31").Value = 3333.9656
$ws.Range("K31").Value = 2401.8125
$ws.Range("L31").Value = 3333.9656
$ws.Range("M31").Value = -2106.8125
$ws.Range("N31").Value = -3923.9656
$ws.Range("H34").Value = 3002.5334
$ws.Range("I34").Value = 2401.8125
$ws.Range("J34").Value = 3333.9656
$ws.Range("K34").Value = 2401.8125
$ws.Range("L34").Value = 3333.9656
$ws.Range("M34").Value = -2199.8125
$ws.Range("N34").Value = -3737.9656
$ws.Range("H58").Value = 21022.076
$ws.Range("I58").Value = 1861.6
$ws.Range("J58").Value = 32997.375
$ws.Range("K58").Value = 1861.6
$ws.Range("L58").Value = 32997.375
$ws.Range("M58").Value = -1658.6
$ws.Range("N58").Value = -33403.375
$ws.Range("H132").Value = 2827.9678
$ws.Range("I132").Value = 2012.125
$ws.Range("K132").Value = 6036.375
$ws.Range("M132").Value = -3506.375
$ws.Range("H134").Value = 1304.7693
$ws.Range("I134").Value = 995.25
$ws.Range("J134").Value = 1800
$ws.Range("K134").Value = 2985.75
$ws.Range("L134").Value = 5400
$ws.Range("M134").Value = -450.75
$ws.Range("N134").Value = -10470
$ws.Range("H136").Value = 21022.076
$ws.Range("I136").Value = 1861.6
$ws.Range("J136").Value = 32997.375
$ws.Range("K136").Value = 5584.799999999999
$ws.Range("L136").Value = 98992.125
$ws.Range("M136").Value = -3034.799999999999
$ws.Range("N136").Value = -104092.125

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 694.03
$ws.Range("I131").Value = 398.625
$ws.Range("J131").Value = 719.7174
$ws.Range("K131").Value = 1195.875
$ws.Range("L131").Value = 2159.1522
$ws.Range("M131").Value = 3844.125
$ws.Range("N131").Value = -12239.1522
$ws.Range("H140").Value = 1482.8889
$ws.Range("I140").Value = 724.3333
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 2172.9999
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = 3007.0001
$ws.Range("N140").Value = -19360

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2848871.8
$ws.Range("I70").Value = 3679.5454
$ws.Range("J70").Value = 5694064
$ws.Range("K70").Value = 3679.5454
$ws.Range("L70").Value = 5694064
$ws.Range("M70").Value = -3409.5454
$ws.Range("N70").Value = -5694604
$ws.Range("H73").Value = 2848871.8
$ws.Range("I73").Value = 3679.5454
$ws.Range("J73").Value = 5694064
$ws.Range("K73").Value = 3679.5454
$ws.Range("L73").Value = 5694064
$ws.Range("M73").Value = -2743.5454
$ws.Range("N73").Value = -5695936
$ws.Range("H80").Value = 3700.7693
$ws.Range("I80").Value = 2855.3333
$ws.Range("J80").Value = 4148.353
$ws.Range("K80").Value = 2855.3333
$ws.Range("L80").Value = 4148.353
$ws.Range("M80").Value = -1857.3333
$ws.Range("N80").Value = -6144.353
$ws.Range("H83").Value = 3700.7693
$ws.Range("I83").Value = 2855.3333
$ws.Range("J83").Value = 4148.353
$ws.Range("K83").Value = 14276.6665
$ws.Range("L83").Value = 20741.765
$ws.Range("M83").Value = -9284.6665
$ws.Range("N83").Value = -30725.765

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3449.8
$ws.Range("J22").Value = 2987.25
$ws.Range("L22").Value = 2987.25
$ws.Range("N22").Value = -3577.25
$ws.Range("H27").Value = 3449.8
$ws.Range("J27").Value = 2987.25
$ws.Range("L27").Value = 2987.25
$ws.Range("N27").Value = -3201.25
$ws.Range("H46").Value = 629.3333
$ws.Range("I46").Value = 580.1
$ws.Range("K46").Value = 580.1
$ws.Range("M46").Value = -392.1
$ws.Range("H100").Value = 2332.8333
$ws.Range("I100").Value = 1000
$ws.Range("J100").Value = 2599.4
$ws.Range("K100").Value = 1000
$ws.Range("L100").Value = 2599.4
$ws.Range("M100").Value = -459
$ws.Range("N100").Value = -3681.4
$ws.Range("H132").Value = 390103.47
$ws.Range("I132").Value = 548803.44
$ws.Range("J132").Value = 2170.3333
$ws.Range("K132").Value = 1646410.32
$ws.Range("L132").Value = 6510.999899999999
$ws.Range("M132").Value = -1643880.32
$ws.Range("N132").Value = -11570.9999
$ws.Range("H136").Value = 2117.2727
$ws.Range("I136").Value = 2048.3333
$ws.Range("K136").Value = 6144.999899999999
$ws.Range("M136").Value = -3594.999899999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1733.3334
$ws.Range("I96").Value = 1600
$ws.Range("K96").Value = 1600
$ws.Range("M96").Value = -227
$ws.Range("H132").Value = 1367.5
$ws.Range("I132").Value = 712.3333
$ws.Range("J132").Value = 2841.625
$ws.Range("K132").Value = 2136.9999
$ws.Range("L132").Value = 8524.875
$ws.Range("M132").Value = 393.0001000000002
$ws.Range("N132").Value = -13584.875
$ws.Range("H136").Value = 22441762
$ws.Range("I136").Value = 26469182
$ws.Range("J136").Value = 3286.4285
$ws.Range("K136").Value = 79407546
$ws.Range("L136").Value = 9859.2855
$ws.Range("M136").Value = -79404996
$ws.Range("N136").Value = -14959.2855
$ws.Range("H137").Value = 48000
$ws.Range("J137").Value = 48000
$ws.Range("L137").Value = 48000
$ws.Range("N137").Value = -58200
